$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.364.33'
$ws.Range("E2").Value = '  -3.85%  '

$ws.Range("D3").Value = '1.974.13'
$ws.Range("E3").Value = '  -5.67%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'238.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.84%  '

$ws.Range("D6").Value = "'0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.09%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = "'55.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.22%  '

$ws.Range("D9").Value = "'58.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.92%  '

$ws.Range("D10").Value = "'0.350"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.56%  '

$ws.Range("D11").Value = "'0.0718"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.95%  '

$ws.Range("E12").Value = '  -5.44%  '

$ws.Range("D13").Value = "'0.878"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.22%  '

$ws.Range("D14").Value = "'14.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.95%  '

$ws.Range("D15").Value = '2.272.35'
$ws.Range("E15").Value = '  -5.22%  '

$ws.Range("D16").Value = "'5.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.52%  '

$ws.Range("D17").Value = '1.976.18'
$ws.Range("E17").Value = '  -5.94%  '

$ws.Range("D18").Value = '35.296.12'
$ws.Range("E18").Value = '  -4.05%  '

$ws.Range("D19").Value = "'16.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.85%  '

$ws.Range("D20").Value = "'69.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.40%  '

$ws.Range("D21").Value = '0.0₃0827'
$ws.Range("E21").Value = '  -6.22%  '

$ws.Range("D22").Value = "'230.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.61%  '

$ws.Range("D23").Value = "'4.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.69%  '

$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("D25").Value = "'2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.84%  '

$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.81%  '

$ws.Range("D27").Value = "'161.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.20%  '

$ws.Range("D28").Value = "'8.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.79%  '

$ws.Range("D29").Value = "'19.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.83%  '

$ws.Range("D30").Value = "'0.118"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.22%  '

$ws.Range("D31").Value = "'1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.77%  '

$ws.Range("D32").Value = "'4.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.68%  '

$ws.Range("D33").Value = "'0.0579"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.25%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = "'0.0880"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.84%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = "'4.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -10.90%  '

$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("E37").Value = '  -2.63%  '

$ws.Range("D38").Value = "'2.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -12.42%  '

$ws.Range("D39").Value = "'4.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.90%  '

$ws.Range("E40").Value = '  -1.37%  '

$ws.Range("D41").Value = "'1.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.17%  '

$ws.Range("D42").Value = "'0.0206"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.75%  '

$ws.Range("D43").Value = "'1.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.16%  '

$ws.Range("D44").Value = '1.365.92'
$ws.Range("E44").Value = '  -2.08%  '

$ws.Range("D45").Value = "'0.0872"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.77%  '

$ws.Range("D46").Value = "'88.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.07%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'7.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.91%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = "'15.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.31%  '

$ws.Range("D49").Value = "'2.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.61%  '

$ws.Range("D50").Value = "'2.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.37%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = "'44.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.47%  '
